$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Outcomes")

# Rows 55-58 describe the Ofsted Leadership Rating counts. They previously
# referenced a "long" dataset with a generic "Count" value column and a
# dimensional filter list(Rating='...'), and have been reworked to use the
# new "ofsted_leadership_data" dataset with a specific *_count value column
# and no dimensional filter.
$ws.Range("F55").Value = "ofsted_leadership_data"
$ws.Range("G55").Value = "outstanding_count"
$ws.Range("I55").Value = "list()"
$ws.Range("J55").Value = "list()"

$ws.Range("F56").Value = "ofsted_leadership_data"
$ws.Range("G56").Value = "good_count"
$ws.Range("I56").Value = "list()"
$ws.Range("J56").Value = "list()"

$ws.Range("F57").Value = "ofsted_leadership_data"
$ws.Range("G57").Value = "requires_improvement_count"
$ws.Range("I57").Value = "list()"
$ws.Range("J57").Value = "list()"

$ws.Range("F58").Value = "ofsted_leadership_data"
$ws.Range("G58").Value = "inadequate_count"
$ws.Range("I58").Value = "list()"
$ws.Range("J58").Value = "list()"

# Update the saved view state of the Outcomes sheet (scrolled back to the
# top-left of the data and the selection moved to I5).
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("I5").Select()
